$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-09-14"

# Update the header label for the "through" column (shared string)
$ws.Range("I1").Value = "2022 (through 09-14)"

# Update September value for the 2022 column
$ws.Range("I10").Value = 65

# Update Total row value for the 2022 column
$ws.Range("I14").Value = 1202
